$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) columns per the latest symbol-list refresh.
$ws.Range("D2").Value = "'329.10"
$ws.Range("E2").Value = "'1.76%"
$ws.Range("D3").Value = "'41.24"
$ws.Range("E3").Value = "'4.92%"
$ws.Range("D4").Value = "'5.642"
$ws.Range("E4").Value = "'-1.44%"
$ws.Range("D5").Value = "'0.08164"
$ws.Range("E5").Value = "'2.08%"
$ws.Range("D6").Value = "'2.016"
$ws.Range("E6").Value = "'2.47%"
$ws.Range("D7").Value = "'8.747"
$ws.Range("E7").Value = "'1.59%"
$ws.Range("D8").Value = "'4.496"
$ws.Range("E8").Value = "'-0.25%"
$ws.Range("E9").Value = "'0.06%"
$ws.Range("D10").Value = "'0.9203"
$ws.Range("E10").Value = "'-0.62%"
$ws.Range("D11").Value = "'0.1272"
$ws.Range("E11").Value = "'3.01%"
$ws.Range("D12").Value = "'0.1956"
$ws.Range("E12").Value = "'0.23%"
$ws.Range("D13").Value = "'0.09350"
$ws.Range("E13").Value = "'1.35%"
$ws.Range("D14").Value = "'0.03804"
$ws.Range("E14").Value = "'4.29%"
$ws.Range("E15").Value = "'0.85%"
$ws.Range("D16").Value = "'0.001299"
$ws.Range("E16").Value = "'1.08%"
$ws.Range("D17").Value = "'0.006097"
$ws.Range("E17").Value = "'-1.35%"
$ws.Range("E19").Value = "'2.84%"
$ws.Range("D21").Value = "'8.152"
$ws.Range("E21").Value = "'-6.39%"
$ws.Range("D22").Value = "'0.1372"
$ws.Range("E22").Value = "'-0.07%"
$ws.Range("E23").Value = "'-1.60%"
$ws.Range("D24").Value = "'0.04397"
$ws.Range("E24").Value = "'-0.22%"
$ws.Range("D25").Value = "'0.001258"
$ws.Range("E25").Value = "'-0.50%"
$ws.Range("D26").Value = "'0.004340"
$ws.Range("E26").Value = "'-5.05%"
$ws.Range("E27").Value = "'4.41%"
$ws.Range("E39").Value = "'11.86%"
$ws.Range("E40").Value = "'1.62%"
$ws.Range("D41").Value = "'0.007449"
$ws.Range("E41").Value = "'-0.08%"
$ws.Range("E42").Value = "'0.94%"
$ws.Range("D43").Value = "'0.008948"
$ws.Range("E43").Value = "'-6.26%"
$ws.Range("D44").Value = "'0.002174"
$ws.Range("E44").Value = "'2.62%"
$ws.Range("D45").Value = "'0.01151"
$ws.Range("E45").Value = "'7.47%"
$ws.Range("D46").Value = "'0.00006593"
$ws.Range("E46").Value = "'-2.75%"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("D48").Value = "'0.003204"
$ws.Range("E48").Value = "'7.80%"
$ws.Range("E49").Value = "'-0.53%"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E51").Value = "'0.03%"
